$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 6422.366246417319
$ws.Range("D2").Value = 10813.16472071209
